{"js": "// Replace the contents of each lattice-multiplication exercise cell in the\n// (single) table with the new problems, while preserving the existing\n// paragraph/run formatting (font size 16pt == w:sz 32) of each cell.\n//\n// Each cell holds one paragraph / one run made of 5 \"lines\" separated by\n// <w:br/> (represented in the Word JS object model as the vertical-tab\n// character \"\\v\"):\n//   1) \"NN x NN\"\n//   2) \"  D1    D2\"\n//   3) \"  ----\"\n//   4) \"C1|    |\"\n//   5) \"C2|    |\"\n\nconst newCellLines = [\n  [\"71 x 70\", \"  7    0\", \"  ----\", \"7|    |\", \"1|    |\"],\n  [\"68 x 99\", \"  9    9\", \"  ----\", \"6|    |\", \"8|    |\"],\n  [\"62 x 74\", \"  7    4\", \"  ----\", \"6|    |\", \"2|    |\"],\n\n  [\"96 x 74\", \"  7    4\", \"  ----\", \"9|    |\", \"6|    |\"],\n  [\"43 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"3|    |\"],\n  [\"19 x 87\", \"  8    7\", \"  ----\", \"1|    |\", \"9|    |\"],\n\n  [\"20 x 29\", \"  2    9\", \"  ----\", \"2|    |\", \"0|    |\"],\n  [\"19 x 88\", \"  8    8\", \"  ----\", \"1|    |\", \"9|    |\"],\n  [\"57 x 77\", \"  7    7\", \"  ----\", \"5|    |\", \"7|    |\"],\n\n  [\"95 x 36\", \"  3    6\", \"  ----\", \"9|    |\", \"5|    |\"],\n  [\"59 x 79\", \"  7    9\", \"  ----\", \"5|    |\", \"9|    |\"],\n  [\"46 x 53\", \"  5    3\", \"  ----\", \"4|    |\", \"6|    |\"],\n\n  [\"38 x 28\", \"  2    8\", \"  ----\", \"3|    |\", \"8|    |\"],\n  [\"98 x 81\", \"  8    1\", \"  ----\", \"9|    |\", \"8|    |\"],\n  [\"88 x 46\", \"  4    6\", \"  ----\", \"8|    |\", \"8|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Collect every cell's first paragraph across the whole table, in\n// row-major (reading) order, so index `i` lines up with `newCellLines[i]`.\nconst cellParagraphs = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < 3; c++) {\n    const cell = table.getCellOrNullObject(r, c);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    cellParagraphs.push(paragraphs);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < newCellLines.length; i++) {\n  const paragraph = cellParagraphs[i].items[0];\n  const newText = newCellLines[i].join(\"\\v\");\n  // Replacing the whole paragraph range's text (rather than clearing the\n  // body first) keeps the existing run formatting (w:rPr / sz) intact.\n  const range = paragraph.getRange(\"Whole\");\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the contents of each lattice-multiplication exercise cell in the\n# (single) table with the new problems, while preserving the existing\n# paragraph/run formatting of each cell (setting Range.Text keeps the\n# run's existing character formatting, e.g. the 16pt / w:sz=32 font).\n#\n# Each cell holds 5 \"lines\" separated by a vertical-tab character\n# (Chr(11)), which Word renders/serializes as <w:br/>:\n#   1) \"NN x NN\"\n#   2) \"  D1    D2\"\n#   3) \"  ----\"\n#   4) \"C1|    |\"\n#   5) \"C2|    |\"\n\n$vt = [char]11\n\n$newCellLines = @(\n  @(\"71 x 70\", \"  7    0\", \"  ----\", \"7|    |\", \"1|    |\"),\n  @(\"68 x 99\", \"  9    9\", \"  ----\", \"6|    |\", \"8|    |\"),\n  @(\"62 x 74\", \"  7    4\", \"  ----\", \"6|    |\", \"2|    |\"),\n\n  @(\"96 x 74\", \"  7    4\", \"  ----\", \"9|    |\", \"6|    |\"),\n  @(\"43 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"3|    |\"),\n  @(\"19 x 87\", \"  8    7\", \"  ----\", \"1|    |\", \"9|    |\"),\n\n  @(\"20 x 29\", \"  2    9\", \"  ----\", \"2|    |\", \"0|    |\"),\n  @(\"19 x 88\", \"  8    8\", \"  ----\", \"1|    |\", \"9|    |\"),\n  @(\"57 x 77\", \"  7    7\", \"  ----\", \"5|    |\", \"7|    |\"),\n\n  @(\"95 x 36\", \"  3    6\", \"  ----\", \"9|    |\", \"5|    |\"),\n  @(\"59 x 79\", \"  7    9\", \"  ----\", \"5|    |\", \"9|    |\"),\n  @(\"46 x 53\", \"  5    3\", \"  ----\", \"4|    |\", \"6|    |\"),\n\n  @(\"38 x 28\", \"  2    8\", \"  ----\", \"3|    |\", \"8|    |\"),\n  @(\"98 x 81\", \"  8    1\", \"  ----\", \"9|    |\", \"8|    |\"),\n  @(\"88 x 46\", \"  4    6\", \"  ----\", \"8|    |\", \"8|    |\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$index = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $lines = $newCellLines[$index]\n    $newText = [string]::Join($vt, $lines)\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newText\n    $index = $index + 1\n  }\n}\n"}
